$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the "Pass/Fail" column (F2:F16) with "PASS" for every test row
$ws.Range("F2:F16").Value = "PASS"

# Move the active selection to I13, matching the saved cursor position
$ws.Range("I13").Select()
